# Apply the "Add files via upload" edit to the risk register sheet.
# Rows 2-9 are rewritten (content shuffled/edited), a brand-new row 9 (R8)
# is appended, some row heights change, and column widths for A/B are
# narrowed. Finally the active selection moves to B13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2: R1 - Faillite de sécurité car serveurs mal configurés ----
$ws.Range("A2").Value = "R1"
$ws.Range("B2").Value = "Faillite de sécurité car serveurs mal configurés"
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 3
$ws.Range("E2").Value = "Mettre à jour la configuration des serveurs pour se conformer à la politique de sécurité de l'entreprise. Mettre en place une équipe en charge de défendre contre les cyber-attaques."
$ws.Range("F2").Value = "Faire remonter une alerte à l'équipe de réponse. Améliorer le procesus de réponse aux incidents et la politique de sécurité en étudiant l'attaque et son impact."
$ws.Range("E2").WrapText = $true
$ws.Range("F2").WrapText = $true
$ws.Rows(2).RowHeight = 119

# ---- Row 3: R2 - Risque antivirus serveurs ----
$ws.Range("A3").Value = "R2"
$ws.Range("B3").Value = "Risque de sécurité car antivirus des serveurs pas mis à jour automatiquement"
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = "Mettre en place des mises à jour automatiques pour l'antivirus. "
$ws.Range("F3").Value = "Faire remonter une alerte à l'équipe de réponse. Améliorer le procesus de réponse aux incidents et la politique de sécurité en étudiant l'attaque et son impact."
$ws.Range("E3").WrapText = $true
$ws.Range("F3").WrapText = $true
$ws.Rows(3).RowHeight = 85

# ---- Row 4: R3 - Risque version Windows serveurs pas patchée ----
$ws.Range("A4").Value = "R3"
$ws.Range("B4").Value = "Risque de sécurité car version de Windows de serveurs pas patchée automatiquement"
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = "Mettre en place des téléchargements automatiques pour les mises à jour Windows. "
$ws.Range("F4").Value = "Faire remonter une alerte à l'équipe de réponse. Améliorer le procesus de réponse aux incidents et la politique de sécurité en étudiant l'attaque et son impact."
$ws.Range("E4").WrapText = $true
$ws.Range("F4").WrapText = $true
$ws.Rows(4).RowHeight = 85

# ---- Row 5: R4 - Congestion du réseau Internet ----
$ws.Range("A5").Value = "R4"
$ws.Range("B5").Value = "Congestion du réseau Internet et raléntissement de l'opérationnel"
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = "Faire des projections d'utilisation du réseau Internet. Acheter de la capacité réseau en conséquence"
$ws.Range("F5").Value = "Achèter de la capacité réseau"
$ws.Range("E5").WrapText = $true
$ws.Range("F5").WrapText = $true
$ws.Rows(5).RowHeight = 68

# ---- Row 6: R5 - Risque accès non autorisé / vol de données ----
$ws.Range("A6").Value = "R5"
$ws.Range("B6").Value = "Risque sur l'accès non autorisé, vol de données et protection de données privées"
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 3
$ws.Range("E6").Value = "Remplacer l'accès TSE administrateur des fournisseurs par un accès avec les droits minimums nécessaire à leurs missions"
$ws.Range("F6").Value = "Revoquer l'accès du  fournisseur en question. Reviser et mettre à jour la politique de droits d'accès de l'entreprise"
$ws.Range("E6").WrapText = $true
$ws.Range("F6").WrapText = $true
$ws.Rows(6).RowHeight = 85

# ---- Row 7 (new): R6 - Violation des données privées (RGPD) ----
$ws.Range("A7").Value = "R6"
$ws.Range("B7").Value = "Violation des données privées (RGPD)"
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 3
$ws.Range("E7").Value = "Encrypter les backups. "
$ws.Range("F7").Value = "Identifier les données qu'on été affectées, mettre en place un plan d'action de réponse et suivre les réglémentations en vigueur du RGPD en ce qui concerne les violations des données"
$ws.Range("E7").WrapText = $true
$ws.Range("F7").WrapText = $true
$ws.Rows(7).RowHeight = 85

# ---- Row 8 (new): R7 - Défaillance ou perte du backup ----
$ws.Range("A8").Value = "R7"
$ws.Range("B8").Value = "Défaillance ou perte du backup"
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 4
$ws.Range("E8").Value = "Procédures de backups et de réstoration doivent être testées régulierement selon la politique de backup. Des supports additionnels doivent être installés et configurés afin d'avoir de backups pour les backups en cas de défaillance de disque."
$ws.Range("F8").Value = "Restaurer le dernier backup valide. Si il n'a pas, restaurer manuellement autant de données que possible. "
$ws.Range("E8").WrapText = $true
$ws.Range("F8").WrapText = $true
$ws.Rows(8).RowHeight = 153

# ---- Row 9 (new): R8 - Défaillance des imprimantes en fin de vie ----
$ws.Range("A9").Value = "R8"
$ws.Range("B9").Value = "Défaillance des imprimantes en fin de vie"
$ws.Range("C9").Value = 2
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = "Remplacer les imprimantes personnelles par des imprimantes réseaux communes (utiliser les imprimantes déjà existantes multi-fonctions pour cela)"
$ws.Range("F9").Value = "Former les utilisateurs à l'utilisation des imprimantes partagées sur le réseau"
$ws.Range("E9").WrapText = $true
$ws.Range("F9").WrapText = $true
$ws.Rows(9).RowHeight = 102

# ---- Column widths: ID column gets narrower, description column too ----
$ws.Columns("A").ColumnWidth = 3.33203125
$ws.Columns("B").ColumnWidth = 73.33203125

# ---- Selection moves to B13 (per the saved view state) ----
$ws.Range("B13").Select()
